$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.565.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.05%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.962.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.17%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.12%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'244.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.21%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -0.48%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'58.69"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.32%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.07%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +3.02%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -7.12%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -0.84%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +0.06%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +0.23%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.250.37"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.23%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'13.78"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.57%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +0.96%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'1.966.63"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.72%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'36.493.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.02%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'69.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.35%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.0₃0857"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.71%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'228.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.82%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E23").Value = "'  -0.13%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.16%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.25%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'9.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.98%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.140"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.17%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'160.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.13%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'19.44"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.19%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +0.94%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -3.03%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -0.23%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -3.82%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -0.02%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -0.15%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +2.13%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +11.53%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -0.17%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'5.77"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -10.06%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -1.92%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'2.90"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.22%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -0.83%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.0212"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.15%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'15.97"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.07%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.366.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.61%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -0.98%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'88.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.74%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -0.91%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +0.15%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'2.141.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.24%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'43.78"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -5.33%  "
$ws.Range("E51").Style = "Normal"
